$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the empty placeholder sheet ("Sheet"); remaining sheets (F, l, h, w) keep their order.
$wb.Worksheets.Item("Sheet").Delete()

# Update column B on sheet "F" with newly computed values (A column is unchanged).
$ws = $wb.Worksheets.Item("F")
$bValues = @(
    [double]"2.279202279202279e-06", [double]"2.329851218740108e-06", [double]"2.380500158277936e-06", [double]"2.431149097815764e-06", [double]"2.481798037353593e-06", [double]"2.532446976891421e-06", [double]"2.58309591642925e-06", [double]"2.633744855967078e-06", [double]"2.684393795504906e-06", [double]"2.735042735042735e-06", [double]"2.785691674580564e-06"
)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $bValues[$i]
}

# Update column B on sheet "l" with newly computed values (A column is unchanged).
$ws = $wb.Worksheets.Item("l")
$bValues = @(
    [double]"1.846153846153847e-06", [double]"1.971986071541628e-06", [double]"2.103409939854384e-06", [double]"2.240547008547009e-06", [double]"2.383518835074391e-06", [double]"2.532446976891422e-06", [double]"2.687452991452993e-06", [double]"2.848658436213993e-06", [double]"3.016184868629315e-06", [double]"3.190153846153848e-06", [double]"3.370686926242484e-06"
)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $bValues[$i]
}

# Update column B on sheet "h" with newly computed values (A column is unchanged).
$ws = $wb.Worksheets.Item("h")
$bValues = @(
    [double]"3.473864165831854e-06", [double]"3.252197256014503e-06", [double]"3.04899561861464e-06", [double]"2.862375869063113e-06", [double]"2.690680516718609e-06", [double]"2.532446976891422e-06", [double]"2.386381347380929e-06", [double]"2.251336140983639e-06", [double]"2.126291315054895e-06", [double]"2.010338058930471e-06", [double]"1.90266489623698e-06"
)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $bValues[$i]
}

# Update column B on sheet "w" with newly computed values (A column is unchanged).
$ws = $wb.Worksheets.Item("w")
$bValues = @(
    [double]"2.813829974323802e-06", [double]"2.752659757490676e-06", [double]"2.694092528607896e-06", [double]"2.637965600928564e-06", [double]"2.584129568256553e-06", [double]"2.532446976891422e-06", [double]"2.482791153815119e-06", [double]"2.435045170087906e-06", [double]"2.389100921595681e-06", [double]"2.344858311936501e-06", [double]"2.302224524446747e-06"
)
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $bValues[$i]
}
